$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.409.59"
$ws.Range("E2").Value = "  -0.41%  "

# Row 3
$ws.Range("D3").Value = "2.179.11"
$ws.Range("E3").Value = "  -1.79%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.43"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.17%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0914"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.14%  "

# Row 12
$ws.Range("E12").Value = "  -0.19%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.92%  "

# Row 14
$ws.Range("D14").Value = "2.503.96"
$ws.Range("E14").Value = "  -2.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.59%  "

# Row 16
$ws.Range("D16").Value = "2.161.39"
$ws.Range("E16").Value = "  -2.87%  "

# Row 17
$ws.Range("E17").Value = "  -4.08%  "

# Row 18
$ws.Range("D18").Value = "42.307.18"
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("E19").Value = "  -3.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

# Row 21
$ws.Range("E21").Value = "  -1.02%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.12%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "

# Row 25
$ws.Range("E25").Value = "  -0.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.19%  "

# Row 27
$ws.Range("E27").Value = "  -0.77%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.03%  "

# Row 29
$ws.Range("E29").Value = "  -2.11%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.80%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "

# Row 33
$ws.Range("E33").Value = "  +2.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.63%  "

# Row 35
$ws.Range("E35").Value = "  -1.47%  "

# Row 36
$ws.Range("E36").Value = "  -0.99%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.96%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0335"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.97%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.71%  "

# Row 40
$ws.Range("E40").Value = "  -4.36%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.196"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.74%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "58.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.92%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.17%  "

# Row 45
$ws.Range("E45").Value = "  +7.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0972"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.68%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.462"
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "  -4.32%  "

# Row 49
$ws.Range("E49").Value = "  -1.94%  "

# Row 50
$ws.Range("E50").Value = "  -0.83%  "

# Row 51
$ws.Range("E51").Value = "  +0.06%  "
